# Applies the "cryptos list" price/volume refresh described by the commit:
# "Updated cryptos list on Tue May 14 10:12:07 UTC 2024 with GitHub Actions"
#
# The Price (column D) and Volume(1h) (column E) cells in this sheet are stored
# as plain text (not numbers/formulas), so every write below is done as text.
# Some Price values (e.g. "7.00", "0.504") look numeric and Excel would silently
# coerce them to a Number (dropping trailing zeros / changing representation) if
# assigned directly, so those cells are forced to Text format first and the cell
# style is reset back to Normal afterwards so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.734.89"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "2.912.09"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "2.910.93"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.40%  "
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000236"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "3.395.08"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").Value = "61.792.69"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.16%  "
$ws.Range("D19").Value = "2.912.18"
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "436.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.659"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.28%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("E29").Value = "  +20.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("E31").Value = "  -2.56%  "
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.108"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.977"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("E43").Value = "  -4.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("D45").Value = "2.697.55"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "342.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.37%  "

Write-Host "Applied cryptos list refresh: updated 82 Price/Volume cells across 50 rows."
